$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BI1").Value = "04-sep"

$values = @(15, 12, 11, 12, 10, 21, 18, 18, 19, 6)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 60).Copy()
    $ws.Cells.Item($row, 61).PasteSpecial(-4122)
    $ws.Cells.Item($row, 61).Value = $values[$i]
}

$ws.Range("BM9").Select()
